$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: Acoustic Guitar
$ws.Range("A22").Value = "acoustic_guitar"
$ws.Range("B22").Value = "Acoustic Guitar"
$ws.Range("C22").Value = "Chitarra Acustica"
$ws.Range("D22").Value = "Akustikgitarre"
$ws.Range("E22").Value = "Akustiskā ģitāra"
$ws.Range("F22").Value = 40
$ws.Range("G22").Value = 88
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = "treble"
$ws.Range("J22").Value = 40
$ws.Range("K22").Value = "NA"

# Row 23: Electric Guitar
$ws.Range("A23").Value = "electric_guitar"
$ws.Range("B23").Value = "Electric Guitar"
$ws.Range("C23").Value = "Chitarra Elettrica"
$ws.Range("D23").Value = "E-Gitarre"
$ws.Range("E23").Value = "Elektriskā ģitāra"
$ws.Range("F23").Value = 40
$ws.Range("G23").Value = 88
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = "treble"
$ws.Range("J23").Value = 40
$ws.Range("K23").Value = "NA"

# Clear the previously-applied cell style (s="3") on A22/A23 so the new
# rows use the default style, matching the target cells (no explicit s attr).
$ws.Range("A22:K23").Style = "Normal"

# Update the selection to match the new active cell/range.
$ws.Range("A22:K23").Select()
